$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 205; this shifts the existing rows 205:260 down to 206:261,
# carrying their original values/formatting with them (matches the diff, where the
# old row 205 data reappears unchanged as row 206, etc., down to the old row 260
# reappearing unchanged as row 261).
$ws.Rows(205).Insert()

# Populate the newly inserted row 205 with the new record's data. Columns that keep
# the same constant value across the whole table (A, B, C, E, F, G, H, I, N, O, Q, R)
# are filled in the same way; D/J/K/L/M/P hold the new record's figures.
$ws.Range("A205").Value = 9
$ws.Range("B205").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C205").Value = "Metropolitana"
$ws.Range("D205").Value = 44798
$ws.Range("E205").Value = 13
$ws.Range("F205").Value = 100112026
$ws.Range("G205").Value = "Haba"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 90
$ws.Range("K205").Value = 12000
$ws.Range("L205").Value = 13000
$ws.Range("M205").Value = 12556
$ws.Range("N205").Value = "$/saco 25 kilos"
$ws.Range("O205").Value = "Provincia de Limarí"
$ws.Range("P205").Value = 502
$ws.Range("Q205").Value = 25
$ws.Range("R205").Value = "Hortaliza"

# Row 205's date cell should carry the same date-time style the other D-column
# cells use (style index 2 / "YYYY-MM-DD HH:MM:SS"); Rows.Insert() already carries
# that style forward from the old row 205, so no further style change is needed.
